$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames ---
$ws.Range("A1").Value = "subid"
$ws.Range("G1").Value = "ethnicity_race"

# --- Fill previously-blank reason_excluded (L) / block_excluded (M) cells with "NA" ---
# L7:L31 were blank -> "NA" (L2:L6 already contain "pilot" and are left untouched)
$lRows = 7..31
foreach ($r in $lRows) {
    $ws.Range("L$r").Value = "NA"
}

# M2:M31 were blank -> "NA", except M13 and M24 which already hold numeric block-excluded values
$mRows = 2..31
foreach ($r in $mRows) {
    if ($r -eq 13 -or $r -eq 24) {
        continue
    }
    $ws.Range("M$r").Value = "NA"
}

# --- Selection moves to L33 ---
$ws.Range("L33").Select()
